$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.470.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.938.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.47%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.515'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.934.06'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.57%  '
$ws.Range('E11').Value = '  -3.62%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('E13').Value = '  -2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.30'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.386.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.426.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.62%  '
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.937.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '445.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.36%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  -1.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.69%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.38'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('E32').Value = '  -4.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.974'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.73'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '45.42'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  -9.35%  '
$ws.Range('E41').Value = '  -6.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.300'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '384.73'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.89%  '
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.703.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.58'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('E50').Value = '  +4.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.25%  '
